# Update "想去人数" (F column) figures across the sheets as per site re-generation.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws2 = $wb.Worksheets.Item("演出")
$ws4 = $wb.Worksheets.Item("全部类型")

# 展览 sheet
$ws1.Cells.Item(2, 6).Value = 829
$ws1.Cells.Item(3, 6).Value = 196
$ws1.Cells.Item(4, 6).Value = 1410
$ws1.Cells.Item(5, 6).Value = 849
$ws1.Cells.Item(6, 6).Value = 483
$ws1.Cells.Item(7, 6).Value = 633
$ws1.Cells.Item(8, 6).Value = 196
$ws1.Cells.Item(9, 6).Value = 9
$ws1.Cells.Item(12, 6).Value = 122
$ws1.Cells.Item(13, 6).Value = 1639
$ws1.Cells.Item(14, 6).Value = 207
$ws1.Cells.Item(15, 6).Value = 35
$ws1.Cells.Item(17, 6).Value = 76
$ws1.Cells.Item(21, 6).Value = 31
$ws1.Cells.Item(23, 6).Value = 733
$ws1.Cells.Item(25, 6).Value = 1473
$ws1.Cells.Item(26, 6).Value = 191

# 演出 sheet
$ws2.Cells.Item(4, 6).Value = 653
$ws2.Cells.Item(5, 6).Value = 202
$ws2.Cells.Item(7, 6).Value = 273

# 全部类型 sheet
$ws4.Cells.Item(3, 6).Value = 829
$ws4.Cells.Item(4, 6).Value = 196
$ws4.Cells.Item(5, 6).Value = 1410
$ws4.Cells.Item(6, 6).Value = 849
$ws4.Cells.Item(9, 6).Value = 483
$ws4.Cells.Item(10, 6).Value = 633
$ws4.Cells.Item(11, 6).Value = 653
$ws4.Cells.Item(12, 6).Value = 196
$ws4.Cells.Item(13, 6).Value = 9
$ws4.Cells.Item(16, 6).Value = 122
$ws4.Cells.Item(17, 6).Value = 1639
$ws4.Cells.Item(18, 6).Value = 202
$ws4.Cells.Item(19, 6).Value = 207
$ws4.Cells.Item(20, 6).Value = 35
$ws4.Cells.Item(22, 6).Value = 76
$ws4.Cells.Item(26, 6).Value = 273
$ws4.Cells.Item(33, 6).Value = 31
$ws4.Cells.Item(35, 6).Value = 733
$ws4.Cells.Item(37, 6).Value = 1473
$ws4.Cells.Item(38, 6).Value = 191
